# Closing and opening balance change
# G2:G11 ("closing_balance" column) -> "2000 : 1"
# H2:H11 ("opening_balance" column) -> "500 : 1"
# Both ranges become wrapped text (matches the style used for the rest of the
# data rows) instead of the old plain numeric 1000 / 500 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$closing = $ws.Range("G2:G11")
$opening = $ws.Range("H2:H11")

$closing.Value = "2000 : 1"
$opening.Value = "500 : 1"

$ws.Range("G2:H11").WrapText = $true

# Leave the selection where the author left it after making the edit.
$null = $ws.Range("G2:H11").Select()
